$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G12 was stored as text ("1234567897"); the edit converts it to a real number.
$ws.Cells.Item(12, 7).Value = 1234567897

# New row 13 - Rohit Sharma
$ws.Cells.Item(13, 1).Value = "Rohit Sharma"
$ws.Cells.Item(13, 2).Value = "rohit@gmail.com"
$ws.Cells.Item(13, 3).Value = "Rohit45"
$ws.Cells.Item(13, 4).Value = "Rohit@987"
$ws.Cells.Item(13, 5).Value = "CEO"
$ws.Cells.Item(13, 6).Value = "Male"
$ws.Cells.Item(13, 7).Value = 7894561237
$ws.Cells.Item(13, 8).Value = "Mumbai"

# New row 14 - Virat kohli
$ws.Cells.Item(14, 1).Value = "Virat kohli"
$ws.Cells.Item(14, 2).Value = "virat@gmail.com"
$ws.Cells.Item(14, 3).Value = "Virat18"
$ws.Cells.Item(14, 4).Value = "Virat@9876"
$ws.Cells.Item(14, 5).Value = "Manager"
$ws.Cells.Item(14, 6).Value = "Male"
$ws.Cells.Item(14, 7).Value = 4561237894
$ws.Cells.Item(14, 8).Value = "Mumbai"

# New row 15 - Dhoni
$ws.Cells.Item(15, 1).Value = "Dhoni"
$ws.Cells.Item(15, 2).Value = "dhoni@gmail.com"
$ws.Cells.Item(15, 3).Value = "Dhoni7"
$ws.Cells.Item(15, 4).Value = "Dhoni@987"
$ws.Cells.Item(15, 5).Value = "Employee"
$ws.Cells.Item(15, 6).Value = "Male"
$ws.Cells.Item(15, 7).Value = 1234567891
$ws.Cells.Item(15, 8).Value = "Csk"

# New row 16 - Abhishek (phone kept as text, matching the source diff)
$ws.Cells.Item(16, 1).Value = "Abhishek"
$ws.Cells.Item(16, 2).Value = "abhi2@gmail.com"
$ws.Cells.Item(16, 3).Value = "Abhi456"
$ws.Cells.Item(16, 4).Value = "Abhi@987"
$ws.Cells.Item(16, 5).Value = "HR"
$ws.Cells.Item(16, 6).Value = "Male"
$ws.Cells.Item(16, 7).NumberFormat = "@"
$ws.Cells.Item(16, 7).Value = "1234567891"
$ws.Cells.Item(16, 8).Value = "abcd"
